# Update the "dSF" (column F) values for several rows to match the
# re-pulled data / corrected mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value  = -3
$ws.Range("F7").Value  = 7
$ws.Range("F8").Value  = 3
$ws.Range("F9").Value  = 3
$ws.Range("F15").Value = 3
